# Update the "Corr/total marks" figures on the marksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking row: correct-answer marks value
$ws.Range("B11").Value = 5

# Total row: total correct marks, and corrected "Corr/total" display string
$ws.Range("B12").Value = 55
$ws.Range("E12").Value = "55/140"
